# Apply ACT_EFF=0 for DC cooling electricity (accounted for through VDA_FLOP)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# The sheet held the ~TFM_TOPINS table; rename it to reflect the broader
# "integration" content now that a ~TFM_INS table is being added below.
$ws.Name = "integration"

# G5:G7 used to repeat the literal values already in F5:F7 ("OUT"/"IN"/"OUT").
# Replace them with a (shared) formula that references column F instead.
$ws.Range("G5:G7").Formula = "=F5"

# --- New ~TFM_INS block --------------------------------------------------
$ws.Range("B9").Value = "~TFM_INS"

# Row 10 is a header row identical to row 4's, copied (with formatting) down.
$ws.Range("B4:I4").Copy()
$ws.Range("B10").PasteSpecial(-4122)
$ws.Range("B10").Value = "TimeSlice"
$ws.Range("C10").Value = "LimType"
$ws.Range("D10").Value = "Attribute"
$ws.Range("E10").Value = "Year"
$ws.Range("F10").Value = "IE"
$ws.Range("G10").Value = "National"
$ws.Range("H10").Value = "Pset_CO"

# Row 11: ACT_EFF=0 for SRVDCE-CS / SRVELC-DC-C (DC cooling electricity),
# since its activity efficiency is already accounted for via VDA_FLOP.
$ws.Range("D11").Value = "ACT_EFF"
$ws.Range("F11").Value = 0
$ws.Range("G11").Formula = "=F11"
$ws.Range("H11").Value = "SRVDCE-CS"
$ws.Range("I11").Value = "SRVELC-DC-C"

# Written last so the shared-string table order matches: ~TFM_INS, ACT_EFF,
# Other_Indexes.
$ws.Range("I10").Value = "Other_Indexes"

# Leave the selection where the workbook was last saved.
$ws.Range("H12").Select()

$wb.Save()
